$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 20, shifting all existing data (rows 20-83)
# down to rows 22-85. This also carries the two oldest rows (old 82/83) down
# to the new rows 84/85, and bumps the sheet dimension automatically.
$ws.Rows("20:21").Insert()

# Populate the newly inserted row 20 with the latest week's data.
$ws.Range("A20").Value = 2
$ws.Range("B20").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C20").Value = "Coquimbo"
$ws.Range("D20").Value = 44497
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 100112043
$ws.Range("G20").Value = "Pepino ensalada"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 600
$ws.Range("K20").Value = 6500
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = 6750
$ws.Range("N20").Value = "$/caja 70 unidades"
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 96
$ws.Range("Q20").Value = 70
$ws.Range("R20").Value = "Hortaliza"

# Populate the newly inserted row 21 with the latest week's data.
$ws.Range("A21").Value = 2
$ws.Range("B21").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44497
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 100112043
$ws.Range("G21").Value = "Pepino ensalada"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Segunda"
$ws.Range("J21").Value = 680
$ws.Range("K21").Value = 4500
$ws.Range("L21").Value = 5000
$ws.Range("M21").Value = 4721
$ws.Range("N21").Value = "$/caja 100 unidades"
$ws.Range("O21").Value = "Provincia de Limarí"
$ws.Range("P21").Value = 47
$ws.Range("Q21").Value = 100
$ws.Range("R21").Value = "Hortaliza"
